$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 1152, pushing the existing rows (old 1152-1247)
# down to (1154-1249). This mirrors the new week's data being prepended to
# the block of "Coliflor" records for this market.
$ws.Rows.Item(1152).EntireRow.Insert()
$ws.Rows.Item(1152).EntireRow.Insert()

# Populate the two freshly-inserted rows with the new weekly observations.
# Row 1152: "Primera" quality
$ws.Cells.Item(1152, 1).Value = 8
$ws.Cells.Item(1152, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1152, 3).Value = "Coquimbo"
$ws.Cells.Item(1152, 4).Value = 45223
$ws.Cells.Item(1152, 5).Value = 4
$ws.Cells.Item(1152, 6).Value = 100112008
$ws.Cells.Item(1152, 7).Value = "Coliflor"
$ws.Cells.Item(1152, 8).Value = "Sin especificar"
$ws.Cells.Item(1152, 9).Value = "Primera"
$ws.Cells.Item(1152, 10).Value = 2000
$ws.Cells.Item(1152, 11).Value = 700
$ws.Cells.Item(1152, 12).Value = 800
$ws.Cells.Item(1152, 13).Value = 750
$ws.Cells.Item(1152, 14).Value = "$/unidad"
$ws.Cells.Item(1152, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1152, 16).Value = 750
$ws.Cells.Item(1152, 17).Value = 1
$ws.Cells.Item(1152, 18).Value = "Hortaliza"

# Row 1153: "Segunda" quality
$ws.Cells.Item(1153, 1).Value = 8
$ws.Cells.Item(1153, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1153, 3).Value = "Coquimbo"
$ws.Cells.Item(1153, 4).Value = 45223
$ws.Cells.Item(1153, 5).Value = 4
$ws.Cells.Item(1153, 6).Value = 100112008
$ws.Cells.Item(1153, 7).Value = "Coliflor"
$ws.Cells.Item(1153, 8).Value = "Sin especificar"
$ws.Cells.Item(1153, 9).Value = "Segunda"
$ws.Cells.Item(1153, 10).Value = 1240
$ws.Cells.Item(1153, 11).Value = 500
$ws.Cells.Item(1153, 12).Value = 600
$ws.Cells.Item(1153, 13).Value = 550
$ws.Cells.Item(1153, 14).Value = "$/unidad"
$ws.Cells.Item(1153, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1153, 16).Value = 550
$ws.Cells.Item(1153, 17).Value = 1
$ws.Cells.Item(1153, 18).Value = "Hortaliza"
